$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("states")

# Update the formula in D2 from 150000-1 to 100000-2 (new population figure of 100000)
$ws.Range("D2").Formula = "=100000-2"
